$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: replace IYR with EEM, update Curr Weight, match the bordered/left-aligned
#     ticker style already used for AGG (A3) / IWM (A5) ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = "EEM"
$ws.Range("C7").Value = 0.2

# --- New row 8: GLD, matching the plain ticker style used for IVV (A4) / EFA (A6) ---
$ws.Range("B6:C6").Copy() | Out-Null
$ws.Range("B8:C8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = "GLD"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.05

# --- Update the sheet selection to cover the ticker column including the new row ---
$ws.Range("A2:A8").Select() | Out-Null
